$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row before existing row 433, pushing rows 433:481 down
# to 434:482 (all of their contents move with them automatically).
$ws.Rows.Item(433).Insert()

# Populate the newly inserted row 433 with the new weekly record. The
# constant columns (A,B,C,E,F,G,H,I,R) match every other row in this block.
$ws.Cells.Item(433, 1).Value = 3
$ws.Cells.Item(433, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(433, 3).Value = "Coquimbo"
$ws.Cells.Item(433, 4).Value = 44946
$ws.Cells.Item(433, 4).NumberFormat = $ws.Cells.Item(434, 4).NumberFormat
$ws.Cells.Item(433, 5).Value = 5
$ws.Cells.Item(433, 6).Value = 100112043
$ws.Cells.Item(433, 7).Value = "Pepino ensalada"
$ws.Cells.Item(433, 8).Value = "Sin especificar"
$ws.Cells.Item(433, 9).Value = "Primera"
$ws.Cells.Item(433, 10).Value = 123
$ws.Cells.Item(433, 11).Value = 15000
$ws.Cells.Item(433, 12).Value = 16000
$ws.Cells.Item(433, 13).Value = 15553
$ws.Cells.Item(433, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(433, 15).Value = "Limache"
$ws.Cells.Item(433, 16).Value = 259
$ws.Cells.Item(433, 17).Value = 60
$ws.Cells.Item(433, 18).Value = "Hortaliza"
